$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (T_s -> T_s Prev, new T_s New column inserted)
$ws.Columns("D:D").Insert()

# Update header row text
$ws.Cells.Item(1, 3).Value2 = "T_s Prev (°C)"
$ws.Cells.Item(1, 4).Value2 = "T_s New (°C)"

# Update data rows (columns C through L, rows 2-19)
$ws.Cells.Item(2, 3).Value2 = 40
$ws.Cells.Item(2, 4).Value2 = 32.57132406623926
$ws.Cells.Item(2, 5).Value2 = 30.5
$ws.Cells.Item(2, 6).Value2 = 30.39742356414612
$ws.Cells.Item(2, 7).Value2 = 7.428675933760736
$ws.Cells.Item(2, 8).Value2 = 2.365886528354384
$ws.Cells.Item(2, 9).Value2 = 5.657153925932682
$ws.Cells.Item(2, 10).Value2 = 4.066524741623234
$ws.Cells.Item(2, 11).Value2 = 5.391634780173662
$ws.Cells.Item(2, 12).Value2 = 4.340845969973043
$ws.Cells.Item(3, 3).Value2 = 32.57132406623926
$ws.Cells.Item(3, 4).Value2 = 33.08246048419874
$ws.Cells.Item(3, 5).Value2 = 30.5
$ws.Cells.Item(3, 6).Value2 = 30.71256737345114
$ws.Cells.Item(3, 7).Value2 = 0.5111364179594773
$ws.Cells.Item(3, 8).Value2 = 1.624696998752291
$ws.Cells.Item(3, 9).Value2 = 5.453972517751855
$ws.Cells.Item(3, 10).Value2 = 4.066524741623234
$ws.Cells.Item(3, 11).Value2 = 5.391634780173662
$ws.Cells.Item(3, 12).Value2 = 4.048610851216464
$ws.Cells.Item(4, 3).Value2 = 33.08246048419874
$ws.Cells.Item(4, 4).Value2 = 33.04534646166731
$ws.Cells.Item(4, 5).Value2 = 30.5
$ws.Cells.Item(4, 6).Value2 = 30.68988423386962
$ws.Cells.Item(4, 7).Value2 = 0.03711402253143348
$ws.Cells.Item(4, 8).Value2 = 1.737128800090907
$ws.Cells.Item(4, 9).Value2 = 5.467743155722625
$ws.Cells.Item(4, 10).Value2 = 4.066524741623234
$ws.Cells.Item(4, 11).Value2 = 5.391634780173662
$ws.Cells.Item(4, 12).Value2 = 4.089581681686743
$ws.Cells.Item(5, 3).Value2 = 33.04534646166731
$ws.Cells.Item(5, 4).Value2 = 33.04801134333719
$ws.Cells.Item(5, 5).Value2 = 30.5
$ws.Cells.Item(5, 6).Value2 = 30.6915088469127
$ws.Cells.Item(5, 7).Value2 = 0.002664881669879549
$ws.Cells.Item(5, 8).Value2 = 1.729916914542366
$ws.Cells.Item(5, 9).Value2 = 5.466742222564275
$ws.Cells.Item(5, 10).Value2 = 4.066524741623234
$ws.Cells.Item(5, 11).Value2 = 5.391634780173662
$ws.Cells.Item(5, 12).Value2 = 4.086934330345509
$ws.Cells.Item(6, 3).Value2 = 33.04801134333719
$ws.Cells.Item(6, 4).Value2 = 33.04781987176892
$ws.Cells.Item(6, 5).Value2 = 30.5
$ws.Cells.Item(6, 6).Value2 = 30.6913921007603
$ws.Cells.Item(6, 7).Value2 = 0.0001914715682644896
$ws.Cells.Item(6, 8).Value2 = 1.730438925463319
$ws.Cells.Item(6, 9).Value2 = 5.466814086714117
$ws.Cells.Item(6, 10).Value2 = 4.066524741623234
$ws.Cells.Item(6, 11).Value2 = 5.391634780173662
$ws.Cells.Item(6, 12).Value2 = 4.08712584987313
$ws.Cells.Item(7, 3).Value2 = 33.04781987176892
$ws.Cells.Item(7, 4).Value2 = 33.04783362832754
$ws.Cells.Item(7, 5).Value2 = 30.5
$ws.Cells.Item(7, 6).Value2 = 30.69140048846765
$ws.Cells.Item(7, 7).Value2 = 0.00001375655861579617
$ws.Cells.Item(7, 8).Value2 = 1.730401440817803
$ws.Cells.Item(7, 9).Value2 = 5.466808923251896
$ws.Cells.Item(7, 10).Value2 = 4.066524741623234
$ws.Cells.Item(7, 11).Value2 = 5.391634780173662
$ws.Cells.Item(7, 12).Value2 = 4.087112096685763
$ws.Cells.Item(8, 3).Value2 = 33.04783362832754
$ws.Cells.Item(8, 4).Value2 = 33.04783362832754
$ws.Cells.Item(8, 5).Value2 = 26.20617762397261
$ws.Cells.Item(8, 6).Value2 = 30.69139988583782
$ws.Cells.Item(8, 7).Value2 = 4.293822376027386
$ws.Cells.Item(8, 8).Value2 = 1.730404134070183
$ws.Cells.Item(8, 9).Value2 = 5.466809294228348
$ws.Cells.Item(8, 10).Value2 = 4.066524741623234
$ws.Cells.Item(8, 11).Value2 = 5.391634780173662
$ws.Cells.Item(8, 12).Value2 = 4.087113084842463
$ws.Cells.Item(9, 3).Value2 = 33.04783362832754
$ws.Cells.Item(9, 4).Value2 = 33.03454381811883
$ws.Cells.Item(9, 5).Value2 = 26.20617762397261
$ws.Cells.Item(9, 6).Value2 = 30.68159114908152
$ws.Cells.Item(9, 7).Value2 = 0.01328981020870401
$ws.Cells.Item(9, 8).Value2 = 2.208418101541191
$ws.Cells.Item(9, 9).Value2 = 5.35273788488761
$ws.Cells.Item(9, 10).Value2 = 4.068000561861331
$ws.Cells.Item(9, 11).Value2 = 5.275361398907608
$ws.Cells.Item(9, 12).Value2 = 4.179155977337425
$ws.Cells.Item(10, 3).Value2 = 33.03454381811883
$ws.Cells.Item(10, 4).Value2 = 33.03543215090832
$ws.Cells.Item(10, 5).Value2 = 26.20617762397261
$ws.Cells.Item(10, 6).Value2 = 30.68212443805691
$ws.Cells.Item(10, 7).Value2 = 0.0008883327894864124
$ws.Cells.Item(10, 8).Value2 = 2.207481154979621
$ws.Cells.Item(10, 9).Value2 = 5.352382849005704
$ws.Cells.Item(10, 10).Value2 = 4.068000561861331
$ws.Cells.Item(10, 11).Value2 = 5.275361398907608
$ws.Cells.Item(10, 12).Value2 = 4.178761256029255
$ws.Cells.Item(11, 3).Value2 = 33.03543215090832
$ws.Cells.Item(11, 4).Value2 = 33.0353727737609
$ws.Cells.Item(11, 5).Value2 = 26.20617762397261
$ws.Cells.Item(11, 6).Value2 = 30.68208879243484
$ws.Cells.Item(11, 7).Value2 = 0.00005937714742287881
$ws.Cells.Item(11, 8).Value2 = 2.207543833292971
$ws.Cells.Item(11, 9).Value2 = 5.352406580075954
$ws.Cells.Item(11, 10).Value2 = 4.068000561861331
$ws.Cells.Item(11, 11).Value2 = 5.275361398907608
$ws.Cells.Item(11, 12).Value2 = 4.178787657361346
$ws.Cells.Item(12, 3).Value2 = 33.0353727737609
$ws.Cells.Item(12, 4).Value2 = 33.0353727737609
$ws.Cells.Item(12, 5).Value2 = 26.3827782945865
$ws.Cells.Item(12, 6).Value2 = 30.68209117503258
$ws.Cells.Item(12, 7).Value2 = 0.1766006706138867
$ws.Cells.Item(12, 8).Value2 = 2.207539644028116
$ws.Cells.Item(12, 9).Value2 = 5.352404993862048
$ws.Cells.Item(12, 10).Value2 = 4.068000561861331
$ws.Cells.Item(12, 11).Value2 = 5.275361398907608
$ws.Cells.Item(12, 12).Value2 = 4.178785892742749
$ws.Cells.Item(13, 3).Value2 = 33.0353727737609
$ws.Cells.Item(13, 4).Value2 = 33.03560093780191
$ws.Cells.Item(13, 5).Value2 = 26.3827782945865
$ws.Cells.Item(13, 6).Value2 = 30.68225665547658
$ws.Cells.Item(13, 7).Value2 = 0.0002281640410117802
$ws.Cells.Item(13, 8).Value2 = 2.194750425321104
$ws.Cells.Item(13, 9).Value2 = 5.3570541434271
$ws.Cells.Item(13, 10).Value2 = 4.067937675067546
$ws.Cells.Item(13, 11).Value2 = 5.280099699425095
$ws.Cells.Item(13, 12).Value2 = 4.177231455218258
$ws.Cells.Item(14, 3).Value2 = 33.03560093780191
$ws.Cells.Item(14, 4).Value2 = 33.03558567300701
$ws.Cells.Item(14, 5).Value2 = 26.3827782945865
$ws.Cells.Item(14, 6).Value2 = 30.68224748975082
$ws.Cells.Item(14, 7).Value2 = 0.00001526479490365773
$ws.Cells.Item(14, 8).Value2 = 2.194766901656863
$ws.Cells.Item(14, 9).Value2 = 5.35706024100541
$ws.Cells.Item(14, 10).Value2 = 4.067937675067546
$ws.Cells.Item(14, 11).Value2 = 5.280099699425095
$ws.Cells.Item(14, 12).Value2 = 4.177238362097799
$ws.Cells.Item(15, 3).Value2 = 33.03558567300701
$ws.Cells.Item(15, 4).Value2 = 33.03558567300701
$ws.Cells.Item(15, 5).Value2 = 26.37818823439179
$ws.Cells.Item(15, 6).Value2 = 30.68224810296311
$ws.Cells.Item(15, 7).Value2 = 0.004590060194715306
$ws.Cells.Item(15, 8).Value2 = 2.194765799360763
$ws.Cells.Item(15, 9).Value2 = 5.357059833060644
$ws.Cells.Item(15, 10).Value2 = 4.067937675067546
$ws.Cells.Item(15, 11).Value2 = 5.280099699425095
$ws.Cells.Item(15, 12).Value2 = 4.177237900014054
$ws.Cells.Item(16, 3).Value2 = 33.03558567300701
$ws.Cells.Item(16, 4).Value2 = 33.03558069513525
$ws.Cells.Item(16, 5).Value2 = 26.37818823439179
$ws.Cells.Item(16, 6).Value2 = 30.68224367495257
$ws.Cells.Item(16, 7).Value2 = 0.000004977871753908403
$ws.Cells.Item(16, 8).Value2 = 2.19510196974334
$ws.Cells.Item(16, 9).Value2 = 5.356938950263645
$ws.Cells.Item(16, 10).Value2 = 4.067939305904179
$ws.Cells.Item(16, 11).Value2 = 5.279976497878478
$ws.Cells.Item(16, 12).Value2 = 4.177279493714362
$ws.Cells.Item(17, 3).Value2 = 33.03558069513525
$ws.Cells.Item(17, 4).Value2 = 33.03558069513525
$ws.Cells.Item(17, 5).Value2 = 26.37830931250751
$ws.Cells.Item(17, 6).Value2 = 30.68224387491598
$ws.Cells.Item(17, 7).Value2 = 0.0001210781157254814
$ws.Cells.Item(17, 8).Value2 = 2.195101610503266
$ws.Cells.Item(17, 9).Value2 = 5.356938817233607
$ws.Cells.Item(17, 10).Value2 = 4.067939305904179
$ws.Cells.Item(17, 11).Value2 = 5.279976497878478
$ws.Cells.Item(17, 12).Value2 = 4.177279343101563
$ws.Cells.Item(18, 3).Value2 = 33.03558069513525
$ws.Cells.Item(18, 4).Value2 = 33.03558118628716
$ws.Cells.Item(18, 5).Value2 = 26.37830931250751
$ws.Cells.Item(18, 6).Value2 = 30.68224399163076
$ws.Cells.Item(18, 7).Value2 = 0.0000004911519084771498
$ws.Cells.Item(18, 8).Value2 = 2.195092745518766
$ws.Cells.Item(18, 9).Value2 = 5.356942005888276
$ws.Cells.Item(18, 10).Value2 = 4.067939262882924
$ws.Cells.Item(18, 11).Value2 = 5.279979747696779
$ws.Cells.Item(18, 12).Value2 = 4.177278246762331
$ws.Cells.Item(19, 3).Value2 = 33.03558118628716
$ws.Cells.Item(19, 4).Value2 = 33.03558118628716
$ws.Cells.Item(19, 5).Value2 = 26.37830629454275
$ws.Cells.Item(19, 6).Value2 = 30.68224397190095
$ws.Cells.Item(19, 7).Value2 = 0.000003017964758100788
$ws.Cells.Item(19, 8).Value2 = 2.195092780964508
$ws.Cells.Item(19, 9).Value2 = 5.356942019013962
$ws.Cells.Item(19, 10).Value2 = 4.067939262882924
$ws.Cells.Item(19, 11).Value2 = 5.279979747696779
$ws.Cells.Item(19, 12).Value2 = 4.177278261623046

Write-Output "Done"
